$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 290.46667
$ws.Range("I33").Value = 287.08334
$ws.Range("J33").Value = 304
$ws.Range("K33").Value = 287.08334
$ws.Range("L33").Value = 304
$ws.Range("M33").Value = -58.08334000000002
$ws.Range("N33").Value = -762
# Row 39
$ws.Range("H39").Value = 677.5909
$ws.Range("I39").Value = 769.9091
$ws.Range("J39").Value = 585.2727
$ws.Range("K39").Value = 2309.7273
$ws.Range("L39").Value = 1755.8181
$ws.Range("M39").Value = -2013.7273
$ws.Range("N39").Value = -2347.8181
# Row 70
$ws.Range("H70").Value = 2757.2856
$ws.Range("I70").Value = 2773.2307
$ws.Range("K70").Value = 8319.6921
$ws.Range("M70").Value = -8049.6921
# Row 73
$ws.Range("H73").Value = 2757.2856
$ws.Range("I73").Value = 2773.2307
$ws.Range("K73").Value = 8319.6921
$ws.Range("M73").Value = -7383.6921
# Row 76
$ws.Range("H76").Value = 4936.75
$ws.Range("I76").Value = 5998
$ws.Range("K76").Value = 5998
$ws.Range("M76").Value = -5683
# Row 79
$ws.Range("H79").Value = 4936.75
$ws.Range("I79").Value = 5998
$ws.Range("K79").Value = 5998
$ws.Range("M79").Value = -4906
# Row 100
$ws.Range("H100").Value = 3963.375
$ws.Range("I100").Value = 2677.5
$ws.Range("J100").Value = 5249.25
$ws.Range("K100").Value = 2677.5
$ws.Range("L100").Value = 5249.25
$ws.Range("M100").Value = -2136.5
$ws.Range("N100").Value = -6331.25
# Row 116
$ws.Range("H116").Value = 3615.0715
$ws.Range("I116").Value = 2789.4443
$ws.Range("K116").Value = 2789.4443
$ws.Range("M116").Value = 652.5556999999999
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents() | Out-Null
$ws.Range("N125").ClearContents() | Out-Null
# Row 132
$ws.Range("H132").Value = 9116.533
$ws.Range("I132").Value = 5527.5557
$ws.Range("J132").Value = 14500
$ws.Range("K132").Value = 16582.6671
$ws.Range("L132").Value = 43500
$ws.Range("M132").Value = -14052.6671
$ws.Range("N132").Value = -48560
# Row 138
$ws.Range("H138").Value = 714713.25
$ws.Range("I138").Value = 1330.1052
$ws.Range("J138").Value = 1037434.2
$ws.Range("K138").Value = 3990.3156
$ws.Range("L138").Value = 3112302.6
$ws.Range("M138").Value = 1149.6844
$ws.Range("N138").Value = -3122582.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1524.238
$ws.Range("I2").Value = 957.82355
$ws.Range("J2").Value = 3931.5
$ws.Range("K2").Value = 957.82355
$ws.Range("L2").Value = 3931.5
$ws.Range("M2").Value = -844.82355
$ws.Range("N2").Value = -4157.5
# Row 45
$ws.Range("H45").Value = 1796.1177
$ws.Range("I45").Value = 1702.2667
$ws.Range("K45").Value = 1702.2667
$ws.Range("M45").Value = -1325.2667
# Row 110
$ws.Range("H110").Value = 1169.9714
$ws.Range("I110").Value = 979.43335
$ws.Range("K110").Value = 979.43335
$ws.Range("M110").Value = 1065.56665
# Row 116
$ws.Range("H116").Value = 1524.238
$ws.Range("I116").Value = 957.82355
$ws.Range("J116").Value = 3931.5
$ws.Range("K116").Value = 957.82355
$ws.Range("L116").Value = 3931.5
$ws.Range("M116").Value = 1336.17645
$ws.Range("N116").Value = -8519.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1524.238
$ws.Range("I3").Value = 957.82355
$ws.Range("J3").Value = 3931.5
$ws.Range("K3").Value = 957.82355
$ws.Range("L3").Value = 3931.5
$ws.Range("M3").Value = -843.82355
$ws.Range("N3").Value = -4159.5
# Row 63
$ws.Range("H63").Value = 30271
$ws.Range("J63").Value = 30271
$ws.Range("L63").Value = 30271
$ws.Range("N63").Value = -31643
# Row 66
$ws.Range("H66").Value = 30271
$ws.Range("J66").Value = 30271
$ws.Range("L66").Value = 90813
$ws.Range("N66").Value = -97677
# Row 99
$ws.Range("H99").Value = 29413068
$ws.Range("I99").Value = 33334604
$ws.Range("J99").Value = 1550
$ws.Range("K99").Value = 33334604
$ws.Range("L99").Value = 1550
$ws.Range("M99").Value = -33333106
$ws.Range("N99").Value = -4546

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1514.8
$ws.Range("I31").Value = 973.63635
$ws.Range("K31").Value = 973.63635
$ws.Range("M31").Value = -678.63635
# Row 34
$ws.Range("H34").Value = 1514.8
$ws.Range("I34").Value = 973.63635
$ws.Range("K34").Value = 973.63635
$ws.Range("M34").Value = -771.63635
# Row 35
$ws.Range("H35").Value = 1196.6666
$ws.Range("I35").Value = 1196.6666
$ws.Range("K35").Value = 1196.6666
$ws.Range("M35").Value = -902.6666
# Row 93
$ws.Range("H93").Value = 37500
$ws.Range("I93").Value = 18750
$ws.Range("K93").Value = 18750
$ws.Range("M93").Value = -16878
# Row 111
$ws.Range("H111").Value = 4702
$ws.Range("J111").Value = 4702
$ws.Range("L111").Value = 4702
$ws.Range("N111").Value = -12882
# Row 116
$ws.Range("H116").Value = 90000
$ws.Range("J116").Value = 90000
$ws.Range("L116").Value = 90000
$ws.Range("N116").Value = -99178

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 421226.12
$ws.Range("I4").Value = 77087.62
$ws.Range("J4").Value = 634264.25
$ws.Range("K4").Value = 231262.86
$ws.Range("L4").Value = 1902792.75
$ws.Range("M4").Value = -231150.86
$ws.Range("N4").Value = -1903016.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2012.3414
$ws.Range("I132").Value = 1674.4814
$ws.Range("J132").Value = 2663.9285
$ws.Range("K132").Value = 5023.4442
$ws.Range("L132").Value = 7991.7855
$ws.Range("M132").Value = -2493.4442
$ws.Range("N132").Value = -13051.7855

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2435.6428
$ws.Range("I40").Value = 2099.4546
$ws.Range("J40").Value = 3668.3333
$ws.Range("K40").Value = 2099.4546
$ws.Range("L40").Value = 3668.3333
$ws.Range("M40").Value = -1963.4546
$ws.Range("N40").Value = -3940.3333
# Row 69
$ws.Range("H69").Value = 57500
$ws.Range("J69").Value = 57500
$ws.Range("L69").Value = 57500
$ws.Range("N69").Value = -59122
# Row 72
$ws.Range("H72").Value = 57500
$ws.Range("J72").Value = 57500
$ws.Range("L72").Value = 172500
$ws.Range("N72").Value = -180612
# Row 74
$ws.Range("H74").Value = 25000
# Row 77
$ws.Range("H77").Value = 25000
# Row 122
$ws.Range("H122").Value = 17859922
$ws.Range("I122").Value = 31252314
$ws.Range("J122").Value = 3399.6667
$ws.Range("K122").Value = 93756942
$ws.Range("L122").Value = 10199.0001
$ws.Range("M122").Value = -93754492
$ws.Range("N122").Value = -15099.0001
# Row 136
$ws.Range("H136").Value = 4199.8
$ws.Range("I136").Value = 4999.6665
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 14998.9995
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -12448.9995
$ws.Range("N136").Value = -14100

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 536.9545
$ws.Range("I113").Value = 276.9375
$ws.Range("J113").Value = 1230.3334
$ws.Range("K113").Value = 830.8125
$ws.Range("L113").Value = 3691.0002
$ws.Range("M113").Value = 1339.1875
$ws.Range("N113").Value = -8031.0002
# Row 132
$ws.Range("H132").Value = 3181.4707
$ws.Range("I132").Value = 3264.6667
$ws.Range("K132").Value = 9794.000100000001
$ws.Range("M132").Value = -7264.000100000001
# Row 136
$ws.Range("H136").Value = 678.8947
$ws.Range("I136").Value = 452.64706
$ws.Range("J136").Value = 2602
$ws.Range("K136").Value = 1357.94118
$ws.Range("L136").Value = 7806
$ws.Range("M136").Value = 1192.05882
$ws.Range("N136").Value = -12906

